$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 22794.77
$ws.Range("J87").Value = 22794.77
$ws.Range("L87").Value = 22794.77
$ws.Range("N87").Value = -25290.77
$ws.Range("H90").Value = 22794.77
$ws.Range("J90").Value = 22794.77
$ws.Range("L90").Value = 68384.31
$ws.Range("N90").Value = -80864.31
$ws.Range("H129").Value = 974.6585
$ws.Range("I129").Value = 382
$ws.Range("J129").Value = 1021.4474
$ws.Range("K129").Value = 1146
$ws.Range("L129").Value = 3064.3422
$ws.Range("M129").Value = 3854
$ws.Range("N129").Value = -13064.3422
$ws.Range("H132").Value = 27032280
$ws.Range("J132").Value = 5883.3335
$ws.Range("L132").Value = 17650.0005
$ws.Range("N132").Value = -22710.0005
$ws.Range("H137").Value = 2084.5
$ws.Range("I137").Value = 1474.4103
$ws.Range("J137").Value = 3119
$ws.Range("K137").Value = 4423.2309
$ws.Range("L137").Value = 9357
$ws.Range("M137").Value = -1873.2309
$ws.Range("N137").Value = -14457
$ws.Range("H141").Value = 2195.6667
$ws.Range("I141").Value = 1851.4348
$ws.Range("J141").Value = 4175
$ws.Range("K141").Value = 5554.3044
$ws.Range("L141").Value = 12525
$ws.Range("M141").Value = -374.3044
$ws.Range("N141").Value = -22885

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7752.7
$ws.Range("I32").Value = 4919.612
$ws.Range("J32").Value = 13504.728
$ws.Range("K32").Value = 4919.612
$ws.Range("L32").Value = 13504.728
$ws.Range("M32").Value = -4632.612
$ws.Range("N32").Value = -14078.728
$ws.Range("H61").Value = 2208.739
$ws.Range("I61").Value = 1849.0714
$ws.Range("J61").Value = 2768.2222
$ws.Range("K61").Value = 1849.0714
$ws.Range("L61").Value = 2768.2222
$ws.Range("M61").Value = -1637.0714
$ws.Range("N61").Value = -3192.2222
$ws.Range("H74").Value = 3427.5117
$ws.Range("I74").Value = 3745.4517
$ws.Range("J74").Value = 2606.1667
$ws.Range("K74").Value = 3745.4517
$ws.Range("L74").Value = 2606.1667
$ws.Range("M74").Value = -2871.4517
$ws.Range("N74").Value = -4354.1667
$ws.Range("H77").Value = 3427.5117
$ws.Range("I77").Value = 3745.4517
$ws.Range("J77").Value = 2606.1667
$ws.Range("K77").Value = 18727.2585
$ws.Range("L77").Value = 13030.8335
$ws.Range("M77").Value = -14359.2585
$ws.Range("N77").Value = -21766.8335
$ws.Range("H132").Value = 2379.348
$ws.Range("I132").Value = 1214.9375
$ws.Range("K132").Value = 3644.8125
$ws.Range("M132").Value = -1114.8125
$ws.Range("H136").Value = 2208.739
$ws.Range("I136").Value = 1849.0714
$ws.Range("J136").Value = 2768.2222
$ws.Range("K136").Value = 5547.2142
$ws.Range("L136").Value = 8304.6666
$ws.Range("M136").Value = -2997.2142
$ws.Range("N136").Value = -13404.6666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 54000
$ws.Range("J63").Value = 54000
$ws.Range("L63").Value = 54000
$ws.Range("N63").Value = -55372
$ws.Range("H66").Value = 54000
$ws.Range("J66").Value = 54000
$ws.Range("L66").Value = 162000
$ws.Range("N66").Value = -168864
$ws.Range("H134").Value = 3029.8293
$ws.Range("I134").Value = 1797.3235
$ws.Range("K134").Value = 5391.970499999999
$ws.Range("M134").Value = -2856.970499999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7688.8887
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("H34").Value = 7688.8887
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("H109").Value = 29998.545
$ws.Range("J109").Value = 29998.545
$ws.Range("L109").Value = 29998.545
$ws.Range("N109").Value = -32078.545
$ws.Range("H122").Value = 3012.5
$ws.Range("I122").Value = 2253.5881
$ws.Range("K122").Value = 6760.7643
$ws.Range("M122").Value = -4310.7643
$ws.Range("H132").Value = 2867.3416
$ws.Range("I132").Value = 2349.5278
$ws.Range("J132").Value = 6595.6
$ws.Range("K132").Value = 7048.5834
$ws.Range("L132").Value = 19786.8
$ws.Range("M132").Value = -4518.5834
$ws.Range("N132").Value = -24846.8
$ws.Range("H134").Value = 4432.405
$ws.Range("I134").Value = 4279.2705
$ws.Range("J134").Value = 5565.6
$ws.Range("K134").Value = 12837.8115
$ws.Range("L134").Value = 16696.8
$ws.Range("M134").Value = -10302.8115
$ws.Range("N134").Value = -21766.8
$ws.Range("M31").ClearContents()
$ws.Range("M34").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 555.3488
$ws.Range("I113").Value = 572.55
$ws.Range("J113").Value = 540.3913
$ws.Range("K113").Value = 1717.65
$ws.Range("L113").Value = 1621.1739
$ws.Range("M113").Value = 452.3500000000001
$ws.Range("N113").Value = -5961.1739
$ws.Range("H132").Value = 2150.3408
$ws.Range("I132").Value = 955.2857
$ws.Range("J132").Value = 2708.0334
$ws.Range("K132").Value = 8597.5713
$ws.Range("L132").Value = 24372.3006
$ws.Range("M132").Value = -6067.5713
$ws.Range("N132").Value = -29432.3006

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 28995.5
$ws.Range("J4").Value = 28995.5
$ws.Range("L4").Value = 28995.5
$ws.Range("N4").Value = -29219.5
$ws.Range("H48").Value = 28999.5
$ws.Range("J48").Value = 28999.5
$ws.Range("L48").Value = 28999.5
$ws.Range("N48").Value = -29969.5
$ws.Range("H102").Value = 2487.2156
$ws.Range("I102").Value = 2110.093
$ws.Range("J102").Value = 4514.25
$ws.Range("K102").Value = 2110.093
$ws.Range("L102").Value = 4514.25
$ws.Range("M102").Value = -488.0929999999998
$ws.Range("N102").Value = -7758.25
$ws.Range("H122").Value = 3174.6924
$ws.Range("I122").Value = 2259
$ws.Range("J122").Value = 4904.3335
$ws.Range("K122").Value = 6777
$ws.Range("L122").Value = 14713.0005
$ws.Range("M122").Value = -4327
$ws.Range("N122").Value = -19613.0005
$ws.Range("H124").Value = 42780
$ws.Range("J124").Value = 42780
$ws.Range("L124").Value = 42780
$ws.Range("N124").Value = -52600
$ws.Range("H132").Value = 2854.875
$ws.Range("I132").Value = 633.9167
$ws.Range("K132").Value = 1901.7501
$ws.Range("M132").Value = 628.2499

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30386
$ws.Range("H100").Value = 1668.5714
$ws.Range("I100").Value = 1136
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1136
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -595
$ws.Range("N100").Value = -4082
$ws.Range("H132").Value = 4409.625
$ws.Range("I132").Value = 973.2
$ws.Range("J132").Value = 16682.572
$ws.Range("K132").Value = 2919.6
$ws.Range("L132").Value = 50047.716
$ws.Range("M132").Value = -389.6000000000004
$ws.Range("N132").Value = -55107.716
$ws.Range("H136").Value = 3582.111
$ws.Range("I136").Value = 1894.4667
$ws.Range("J136").Value = 5691.6665
$ws.Range("K136").Value = 5683.4001
$ws.Range("L136").Value = 17074.9995
$ws.Range("M136").Value = -3133.4001
$ws.Range("N136").Value = -22174.9995

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 40786.25
$ws.Range("J125").Value = 40786.25
$ws.Range("L125").Value = 40786.25
$ws.Range("N125").Value = -50626.25
$ws.Range("H127").Value = 39429.75
$ws.Range("J127").Value = 39429.75
$ws.Range("L127").Value = 39429.75
$ws.Range("N127").Value = -49349.75
$ws.Range("H136").Value = 3156.0938
$ws.Range("I136").Value = 854.2778
$ws.Range("J136").Value = 6115.5713
$ws.Range("K136").Value = 2562.8334
$ws.Range("L136").Value = 18346.7139
$ws.Range("M136").Value = -12.83339999999998
$ws.Range("N136").Value = -23446.7139
